# Spanish translations for "Email 5-1 [TEMPLATE] Partner email – invite revoked"
$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $d.Content.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 2) | Out-Null
}

# Language switcher line
Replace-Text "English" "Inglés"
Replace-Text " / Portuguese / French / Thai / Vietnamese / Spanish" " / Portugués / Francés / Tailandés / Vietnamita / Español"

# Brief / Target audience table labels
Replace-Text "Brief" "Breve"
Replace-Text "An email sent to partners in the target country who RSVPed yes but didn’t submit their documents by the deadline. We will be revoking their invites. It will be sent via customer.io" "An email sent to partners in the target country who RSVPed yes but didn’t submit their documents by the deadline. We will be revoking their invites. Se enviará a través de customer.io"
Replace-Text "Target audience" "Público objetivo"

# Heading
Replace-Text "We didn’t receive your documents on time" "No hemos recibido tus documentos a tiempo"

# Greeting
Replace-Text "Hi " "Hola "

# Body paragraph
Replace-Text "We didn’t receive your documents by the deadline (" "No hemos recibido tus documentos dentro del plazo ("
Replace-Text "). Unfortunately, we’re unable to proceed with your registration for the " "). Lamentablemente, no podemos proceder con tu inscripción para el "

# Contact paragraphs
Replace-Text "If you have any questions, please contact us via " "Si tienes alguna pregunta, entra en contacto con nosotros por "
Replace-Text " or " " o "
Replace-Text "If you have any questions, please contact your country manager, " "Si tienes alguna pregunta, entra en contacto con el gestor de tu país "
Replace-Text ", at " ", en "

# Comments ("choose either one" appears twice, identical text in both comments)
foreach ($c in $d.Comments) {
    $c.Range.Text = "elija uno de los dos"
}
